$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.434.32"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "1.797.35"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.33"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3805"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3457"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.86"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07522"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.14"
$ws.Range("E13").Value = "  +8.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.489"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "1.797.66"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.088"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001107"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06662"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.99"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.548"
$ws.Range("E21").Value = "  +5.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.43"
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("D23").Value = "27.431.32"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.57"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.425"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.568"
$ws.Range("E26").Value = "  +6.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.494"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.50"
$ws.Range("E28").Value = "  +9.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.70"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "2.001.34"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.99"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.043"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.121"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08712"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.31"
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.653"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.466"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6929"
$ws.Range("E38").Value = "  +10.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.896"
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06395"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2203"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02348"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.273"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.54"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6473"
$ws.Range("E45").Value = "  +6.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.873"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.138"
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07205"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.90"
$ws.Range("E51").Value = "  +2.45%  "
